$d = $word.ActiveDocument

# Locate the insertion point: the paragraph that carries the "_GoBack"
# bookmark (the paragraph right after the empty sz=19 spacer paragraph at
# the end of the "delete" Q&A section). New content goes immediately
# before it, matching the XML diff.
$bm = $d.Bookmarks.Item("_GoBack")
$anchorPara = $bm.Range.Paragraphs.Item(1)
$insertPoint = $d.Range($anchorPara.Range.Start, $anchorPara.Range.Start)

$frag = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="19"/></w:rPr><w:t>3.Using</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> pointer</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>Why should we use delete?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>We use delete operator for delete link between pointer and value that pointer is pointting.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">When we use </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>delete ?</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">We use delete when we want give back memmory on Heap for operating systeam and after we use </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/><w:u w:val="single"/></w:rPr><w:t>new</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> operator </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">Difference between delete and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>delete[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>]? Write a demo.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>Delete[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>] for array</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>Delete not array.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/><w:u w:val="single"/></w:rPr><w:t>DEMO</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>#include&lt;iostream&gt;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">Using namespace </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>std</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">Void </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>main(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>){</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">Int *array = new </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>int[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>5];</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">Delete </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>array[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>0];</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>Int *p = new int;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/></w:rPr><w:t>Delete;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="23"/><w:szCs w:val="19"/><w:u w:val="single"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr><w:t>4. Given below code</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">What </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr><w:t>are the result</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr><w:t>?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr><w:t>1. *p1 == 10</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr><w:t>2. *p2 == 0x100</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr><w:t>3. *(*p2) ==</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/><w:color w:val="000000"/><w:sz w:val="25"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> 10</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($frag)

# InsertXML merges the *last* paragraph of the inserted fragment into the
# paragraph that followed the insertion point (rather than creating a
# separate empty paragraph before it), so the fragment above carries one
# extra trailing empty paragraph, which ends up duplicated just before the
# "_GoBack" paragraph. Remove that leftover empty paragraph so the
# surrounding structure is restored exactly.
$bm2 = $d.Bookmarks.Item("_GoBack")
$anchorPara2 = $bm2.Range.Paragraphs.Item(1)
$leftoverIndex = $anchorPara2.Index - 1
$leftover = $d.Paragraphs.Item($leftoverIndex)
$leftover.Range.Delete()

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
